$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "Players" sheet: row 8 (Slashy) - increase L8 (Crit column)
#    from 4 to 30 and add two new weapon slots (P8, Q8).
# ---------------------------------------------------------------------------
$players = $wb.Worksheets.Item("Players")
$players.Range("L8").Value = 30
$players.Range("P8").Value = "Vampic Sword"
$players.Range("Q8").Value = "Seven Sword"
$players.Range("Q9").Select()

# ---------------------------------------------------------------------------
# 2. Add three new lookup sheets (Drain Reversal, Drain Success, Multi-hit)
#    using "Reflect" as the template, inserted right before "Reflect".
# ---------------------------------------------------------------------------

# --- Drain Reversal ---
$wb.Worksheets.Item("Reflect").Copy($wb.Worksheets.Item("Reflect"))
$drainReversal = $wb.Worksheets.Item("Reflect (2)")
$drainReversal.Name = "Drain Reversal"
$drainReversal.Range("A2").Value = "Slashy"
$drainReversal.Range("K2").Value = "Vampic Sword"
$drainReversal.Range("L2").Value = "Skelton"
$drainReversal.Range("A3").Value = "Skelton"
$drainReversal.Range("C3").Value = 6
$drainReversal.Range("K3").Select()

# --- Drain Success ---
$wb.Worksheets.Item("Reflect").Copy($wb.Worksheets.Item("Reflect"))
$drainSuccess = $wb.Worksheets.Item("Reflect (2)")
$drainSuccess.Name = "Drain Success"
$drainSuccess.Range("A2").Value = "Slashy"
$drainSuccess.Range("K2").Value = "Vampic Sword"
$drainSuccess.Range("L2").Value = "Moth"
$drainSuccess.Range("A3").Value = "Moth"
$drainSuccess.Range("C3").Value = 8
$drainSuccess.Range("K3").Select()

# --- Multi-hit ---
$wb.Worksheets.Item("Reflect").Copy($wb.Worksheets.Item("Reflect"))
$multiHit = $wb.Worksheets.Item("Reflect (2)")
$multiHit.Name = "Multi-hit"
$multiHit.Range("A2").Value = "Slashy"
$multiHit.Range("K2").Value = "Seven Sword"
$multiHit.Range("L2").Value = "Crab"
$multiHit.Range("A3").Value = "Crab"
$multiHit.Range("C3").Value = 6
$multiHit.Range("L3").Select()

# ---------------------------------------------------------------------------
# 3. Make "Multi-hit" the active sheet/tab, matching the saved workbook state.
# ---------------------------------------------------------------------------
$multiHit.Activate()
